# Updated cryptos list on Mon Jul 29 20:28:27 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for each coin row.
# Column D values are numeric-looking text (e.g. "67.221.15", "186.39")
# that must stay stored as literal text, exactly like the source data, so
# every Price assignment is prefixed with a leading apostrophe - this is
# the standard Excel idiom that forces text entry instead of having the
# value auto-parsed/reformatted as a number (which would also silently
# drop meaningful trailing zeros, e.g. "446.30" -> 446.3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.221.15"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "'3.322.91"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'186.39"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'577.78"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("D7").Value = "'0.607"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.130"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "'3.890.68"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'27.43"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "'67.501.37"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'3.333.35"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "'446.30"
$ws.Range("E18").Value = "  +7.22%  "
$ws.Range("D19").Value = "'5.69"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D22").Value = "'74.14"
$ws.Range("E22").Value = "  +3.70%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'3.459.15"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "'5.34"
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'6.85"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.24"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  +5.68%  "
$ws.Range("D37").Value = "'162.99"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'27.27"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "'2.781.66"
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").Value = "'4.49"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "'6.25"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0672"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'40.11"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.41"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'325.84"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.989"
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'31.20"
$ws.Range("E51").Value = "  +1.51%  "
